$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 662, pushing the existing data
# (old rows 662:676) down to become rows 668:682.
$ws.Rows("662:667").Insert()

# Static columns shared by every Durazno row in this block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria = "Durazno"

# New rows to fill in, in order, for rows 662..667.
$newRows = @(
    @{ Row=662; D=44595; K="Andross"; L="Especial";                 M=200; N=8000;  O=8000;  P=8000;  Q="$/bandeja 8 kilos empedrada"; R="Provincia de San Felipe de Aconcagua"; S=1000; T=8 },
    @{ Row=663; D=44595; K="Andross"; L="Extra (doble especial)";   M=220; N=10000; O=10000; P=10000; Q="$/bandeja 8 kilos empedrada"; R="Provincia de San Felipe de Aconcagua"; S=1250; T=8 },
    @{ Row=664; D=44595; K="Andross"; L="Primera";                  M=280; N=6000;  O=6000;  P=6000;  Q="$/bandeja 8 kilos empedrada"; R="Provincia de San Felipe de Aconcagua"; S=750;  T=8 },
    @{ Row=665; D=44595; K="Carson";  L="Especial";                 M=180; N=18000; O=18000; P=18000; Q="$/caja 15 kilos granel";       R="Región de O'Higgins";                 S=1200; T=15 },
    @{ Row=666; D=44595; K="Carson";  L="Primera";                  M=200; N=15000; O=15000; P=15000; Q="$/caja 15 kilos granel";       R="Región de O'Higgins";                 S=1000; T=15 },
    @{ Row=667; D=44595; K="Carson";  L="Segunda";                  M=150; N=12000; O=12000; P=12000; Q="$/caja 15 kilos granel";       R="Región de O'Higgins";                 S=800;  T=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
